$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EmulatorData")

# Add new row 14 - AutoPay test case, mirroring row 13's layout but with the
# next sequential ID (13) and Action set to AutoPay.
$ws.Range("A14").Value = "No Emulator Data"
$ws.Range("C14").Value = "13"
$ws.Range("D14").Value = "2.5"
$ws.Range("E14").Value = "10.50"
$ws.Range("G14").Value = "AutoPay"
$ws.Range("H14").Value = "en_US"
$ws.Range("W14").Value = "udf data 4"
$ws.Range("X14").Value = "udf data 5"
$ws.Range("AB14").Value = "udf data 9"
$ws.Range("AC14").Value = "udf data 10"

$ws.Range("AH14").Select()
